# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (column E) and "Valor Mora" (column F) history for
# rows 16-85 is refreshed: the periods now run in ascending chronological
# order (1608 .. 2205 instead of 2205 .. 1608) and the historical "Valor
# Mora" amounts were corrected to match the new period ordering. Only the
# cell values change - existing number formats / borders stay as they are.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$periods = @("1608","1609","1610","1611","1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112","2201","2202","2203","2204","2205")

$amounts = @(27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,23958)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $amounts[$i]
}
